$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts MSSV -> C, Lời nhắn -> D)
$ws.Range("B1").EntireColumn.Insert()

# Set the new header in B1
$ws.Range("B1").Value = "Ngày"

# Update selection to B1
$ws.Range("B1").Select()
